# Resort the worksheets: put "总计" before "2022-Q2" (swap their tab order),
# keeping "2022-Q2" as the active/selected sheet.
$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("总计")
$detailSheet  = $wb.Worksheets.Item("2022-Q2")

# Move "总计" so it becomes the first tab (i.e. move it before "2022-Q2").
$summarySheet.Move($detailSheet)

# Re-fetch "2022-Q2" by name (worksheet handles are position-bound) and keep
# it as the selected/active sheet after reordering.
$wb.Worksheets.Item("2022-Q2").Select()
